$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look like plain numbers need to be
# forced to Text format first, since these columns store text-formatted
# price strings (e.g. "609.69"), not numeric cells.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values.
$ws.Range("D2").Value = "68.432.96"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "3.803.26"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").Value = "609.69"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("D6").Value = "163.18"
$ws.Range("E6").Value = "  -1.64%  "
$ws.Range("D7").Value = "3.800.01"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").Value = "0.516"
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("E10").Value = "  -0.05%  "
$ws.Range("D11").Value = "6.88"
$ws.Range("E11").Value = "  +8.91%  "
$ws.Range("D12").Value = "0.449"
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("E13").Value = "  -1.63%  "
$ws.Range("D14").Value = "35.01"
$ws.Range("E14").Value = "  -2.72%  "
$ws.Range("D15").Value = "4.441.30"
$ws.Range("E15").Value = "  -0.07%  "
$ws.Range("D16").Value = "3.857.06"
$ws.Range("E16").Value = "  +2.24%  "
$ws.Range("D17").Value = "68.380.67"
$ws.Range("E17").Value = "  +0.83%  "
$ws.Range("D18").Value = "18.00"
$ws.Range("E18").Value = "  -2.21%  "
$ws.Range("E19").Value = "  +0.38%  "
$ws.Range("D20").Value = "7.05"
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("D21").Value = "461.41"
$ws.Range("E21").Value = "  -0.44%  "
$ws.Range("E22").Value = "  -2.78%  "
$ws.Range("D23").Value = "0.696"
$ws.Range("E23").Value = "  -0.77%  "
$ws.Range("D24").Value = "0.0000146"
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").Value = "83.31"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "11.96"
$ws.Range("E27").Value = "  -0.53%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").Value = "3.949.57"
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("E31").Value = "  -5.91%  "
$ws.Range("D32").Value = "2.19"
$ws.Range("E32").Value = "  -1.45%  "
$ws.Range("D33").Value = "7.20"
$ws.Range("E33").Value = "  -2.64%  "
$ws.Range("D34").Value = "28.86"
$ws.Range("E34").Value = "  -2.07%  "
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("D36").Value = "9.05"
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("D37").Value = "0.100"
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("E38").Value = "  +5.84%  "
$ws.Range("D39").Value = "5.85"
$ws.Range("E39").Value = "  +0.54%  "
$ws.Range("D40").Value = "0.978"
$ws.Range("E40").Value = "  -1.94%  "
$ws.Range("D41").Value = "3.16"
$ws.Range("E41").Value = "  -1.04%  "
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("B44").Value = "ONDO"
$ws.Range("C44").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D44").Value = "1.41"
$ws.Range("E44").Value = "  +2.62%  "
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").Value = "153.02"
$ws.Range("E45").Value = "  +1.27%  "
$ws.Range("B46").Value = "Arweave"
$ws.Range("C46").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D46").Value = "43.05"
$ws.Range("E46").Value = "  -2.93%  "
$ws.Range("B47").Value = "TheGraph"
$ws.Range("C47").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D47").Value = "0.296"
$ws.Range("E47").Value = "  -1.28%  "
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").Value = "46.87"
$ws.Range("E48").Value = "  -1.72%  "
$ws.Range("D49").Value = "8.34"
$ws.Range("E49").Value = "  -0.24%  "
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("D51").Value = "378.67"
$ws.Range("E51").Value = "  -2.60%  "
